# Fix formatting issue in grants section: remove the asterisks that were
# wrapping "PI" / "Co-PI" markers (e.g. "*PI*" -> "PI", "*Co-PI*" -> "Co-PI"),
# including fixing a missing comma on the Georgia Ornithological Society entry.
# Also corrects the lab-transfer year noted on the education tab (2020 -> 2022).

$wb = $excel.ActiveWorkbook

# --- education sheet: fix advisor note year ---
$eduWs = $wb.Worksheets.Item("education")
$eduWs.Range("E4").Value = "Advisor: Dr. Nathan Senner (Lab transferred institutions in Fall 2022)"

# --- grants sheet: drop asterisks around PI / Co-PI ---
$grantsWs = $wb.Worksheets.Item("grants")

# Update the short role column (E) first, then the funder/amount column (B),
# mirroring the order the text was originally retyped in.
$grantsWs.Range("E4").Value = "PI"
$grantsWs.Range("E5").Value = "Co-PI"

$grantsWs.Range("B4").Value = "American Ornithological Society (`$2,420, PI)"
$grantsWs.Range("B5").Value = "National Fish and Wildlife Foundation (`$90,000, Co-PI)"
$grantsWs.Range("B6").Value = "Association of Field Ornithologists (`$1,476, PI)"
$grantsWs.Range("E6").Value = "PI"

$grantsWs.Range("B7").Value = "Georgia Ornithological Society (`$7,500, PI)"
$grantsWs.Range("E7").Value = "PI"

$grantsWs.Range("B9").Value = "National Fish and Wildlife Foundation (`$84,418 Co-PI)"
$grantsWs.Range("E9").Value = "Co-PI"

$grantsWs.Range("B10").Value = "Carolina Bird Club (`$7,500, PI)"
$grantsWs.Range("E10").Value = "PI"

# --- restore view state: grants tab active, selection on B11 ---
$grantsWs.Range("B11").Select()
$grantsWs.Activate()
